$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166428089141846
$ws.Range("B1").Value = 2.437504291534424
$ws.Range("D1").Value = 2.368390798568726
$ws.Range("E1").Value = 1.234018802642822
